# Get Driver Test cases updated
$wb = $excel.ActiveWorkbook

$wsDriverPost = $wb.Worksheets.Item("driver_Post")
$wsValidation = $wb.Worksheets.Item("driver_Post_ValidationData")

# Update the test data on the driver_Post sheet (new mobile number / email),
# writing the mobile number first so the new shared strings land in the same
# order as the authored edit (mobile number then email).
$wsDriverPost.Range("B5").Value = "7755663245"
$wsDriverPost.Range("B4").Value = "driver.test_45@gmail.com"

# Move the selection/active cell on the validation sheet away from B5, and
# clear it as the tab-selected sheet.
$wsValidation.Range("L6").Select() | Out-Null

# Make driver_Post the active/selected sheet & tab, with B4 as the active cell.
$wsDriverPost.Activate() | Out-Null
$wsDriverPost.Range("B4").Select() | Out-Null
